# Insert a new row 23 (pushing current rows 23:54 down to 24:55),
# then populate the new row 23 with the latest weekly price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 23:54 down by one row.
$ws.Rows(23).Insert()

# Fill the new row 23 with the new weekly record.
$ws.Range("A23").Value = 11
$ws.Range("B23").Value = "Vega Monumental Concepción"
$ws.Range("C23").Value = "Bíobío"
$ws.Range("D23").Value = 44580
$ws.Range("E23").Value = 8
$ws.Range("F23").Value = 100112012
$ws.Range("G23").Value = "Espinaca"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 10000
$ws.Range("L23").Value = 11000
$ws.Range("M23").Value = 10500
$ws.Range("N23").Value = "$/cuna 10 kilos"
$ws.Range("O23").Value = "Región Metropolitana"
$ws.Range("P23").Value = 1050
$ws.Range("Q23").Value = 10
$ws.Range("R23").Value = "Hortaliza"

# Match the date column's number format/style used by the rest of column D.
$ws.Range("D23").NumberFormat = $ws.Range("D24").NumberFormat
